$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve a plain (unbordered) style reference cell to restore style on text-forced numeric-looking cells
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "68.506.35"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "2.693.94"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.14"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.67"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  +1.93%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "2.692.77"
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.359"
$ws.Range("D13").Style = $plainStyle
$ws.Range("E13").Value = "  +2.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.22"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").Value = "3.185.86"
$ws.Range("E15").Value = "  +1.87%  "
$ws.Range("E16").Value = "  -1.13%  "
$ws.Range("D17").Value = "68.438.94"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "2.695.63"
$ws.Range("E18").Value = "  +2.43%  "
$ws.Range("E19").Value = "  +4.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "364.67"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.62"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  +3.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.52"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = "  +2.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.88"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  +1.90%  "
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.33"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = "  -1.52%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.86"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "582.49"
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = "  +4.77%  "
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("E32").Value = "  +1.47%  "
$ws.Range("E33").Value = "  +2.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.94"
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = "  +5.26%  "
$ws.Range("E35").Value = "  +3.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.63"
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = "  +5.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.73"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("E41").Value = "  +2.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.36"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = "  +0.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.66"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  +1.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.85"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0316"
$ws.Range("E46").Value = "  -5.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "157.33"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.76"
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = "  +4.51%  "
$ws.Range("E50").Value = "  +6.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.97"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  -0.20%  "
